$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2..14 in column D hold AA-sequence strings that currently start with
# the "TESTTTTT" placeholder prefix. The commit replaces that prefix with
# "TTEESSTT" while leaving the rest of each sequence untouched.
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $old = $cell.Value2
    if ($old.Length -ge 8) {
        $new = "TTEESSTT" + $old.Substring(8)
        $cell.Value = $new
    }
}

# The selection/view state for column D changed: the whole column D is now
# selected (previously only D14 was selected).
$ws.Columns.Item(4).Select()
